# Updates cryptos list values (prices/volume deltas) per the Sun Jan 14
# 2024 03:22:30 UTC GitHub Actions data refresh. Several rows also swap
# which coin occupies them, reflecting a re-sort of the ranking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores plain-looking numbers (e.g. "302.83") as TEXT
# in this sheet (t="inlineStr" in the OOXML), not as numeric cells.
# Assigning such a string straight to .Value makes Excel's COM layer
# auto-coerce it to a float (introducing binary rounding noise and
# flipping the cell's type), so for any replacement price that parses as
# a plain number we instead prefix it with an apostrophe (Excel's
# "treat as text" marker), which keeps the value/style and write it back
# as text, then restore the cell's original (default/"Normal") style so
# no stray number-format is left behind.
function Set-TextValue([string]$cellRef, [string]$text) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($cellRef).Style = 'Normal'
}

$ws.Range('D2').Value = '42.749.09'
$ws.Range('E2').Value = '  +0.52%  '
$ws.Range('D3').Value = '2.561.97'
$ws.Range('E3').Value = '  +0.73%  '
$ws.Range('E4').Value = '  -0.02%  '
Set-TextValue 'D5' '302.83'
$ws.Range('E5').Value = '  +2.57%  '
Set-TextValue 'D6' '97.31'
$ws.Range('E6').Value = '  +4.99%  '
$ws.Range('E7').Value = '  +0.84%  '
Set-TextValue 'D9' '0.547'
$ws.Range('E9').Value = '  +0.55%  '
Set-TextValue 'D10' '36.48'
$ws.Range('E10').Value = '  +3.51%  '
Set-TextValue 'D11' '0.0810'
$ws.Range('E11').Value = '  +1.19%  '
$ws.Range('E12').Value = '  +9.83%  '
Set-TextValue 'D13' '7.69'
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('D14').Value = '2.555.18'
$ws.Range('E14').Value = '  +0.89%  '
Set-TextValue 'D15' '0.882'
$ws.Range('E15').Value = '  +2.81%  '
Set-TextValue 'D16' '14.54'
$ws.Range('E16').Value = '  +3.90%  '
$ws.Range('D17').Value = '42.822.22'
$ws.Range('E17').Value = '  +0.64%  '
Set-TextValue 'D18' '13.59'
$ws.Range('E18').Value = '  +8.66%  '
$ws.Range('D19').Value = '0.0₃0989'
$ws.Range('E19').Value = '  +2.94%  '
$ws.Range('E20').Value = '  +0.65%  '
Set-TextValue 'D21' '71.67'
$ws.Range('E21').Value = '  -0.65%  '
Set-TextValue 'D22' '257.27'
$ws.Range('E22').Value = '  +0.01%  '
Set-TextValue 'D23' '2.96'
$ws.Range('E23').Value = '  +3.39%  '
Set-TextValue 'D24' '2.10'
$ws.Range('E24').Value = '  +0.00%  '
Set-TextValue 'D25' '28.21'
$ws.Range('E25').Value = '  -4.16%  '
$ws.Range('E26').Value = '  -0.06%  '
Set-TextValue 'D27' '39.22'
$ws.Range('E27').Value = '  +9.80%  '
Set-TextValue 'D28' '10.11'
$ws.Range('E28').Value = '  +2.13%  '
Set-TextValue 'D29' '2.11'
$ws.Range('E29').Value = '  +0.02%  '
Set-TextValue 'D30' '6.01'
$ws.Range('E30').Value = '  +2.47%  '
Set-TextValue 'D31' '156.57'
$ws.Range('E31').Value = '  +4.31%  '
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D32' '2.17'
$ws.Range('E32').Value = '  +1.40%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D33' '2.75'
$ws.Range('E33').Value = '  +1.08%  '
Set-TextValue 'D34' '27.16'
$ws.Range('E34').Value = '  +12.05%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D35' '3.36'
$ws.Range('E35').Value = '  -0.32%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D36' '0.0802'
$ws.Range('E36').Value = '  +2.13%  '
$ws.Range('E37').Value = '  +2.51%  '
Set-TextValue 'D38' '18.19'
$ws.Range('E38').Value = '  +16.89%  '
$ws.Range('E39').Value = '  +0.90%  '
$ws.Range('E40').Value = '  +2.38%  '
Set-TextValue 'D41' '2.04'
$ws.Range('E41').Value = '  +27.82%  '
$ws.Range('B42').Value = 'NEARProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D42' '3.37'
$ws.Range('E42').Value = '  -0.30%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D43' '0.0305'
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.061.92'
$ws.Range('E44').Value = '  +0.14%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D45' '0.999'
$ws.Range('E45').Value = '  +0.05%  '
Set-TextValue 'D46' '88.45'
$ws.Range('E46').Value = '  +5.12%  '
$ws.Range('E47').Value = '  +6.96%  '
Set-TextValue 'D48' '77.03'
$ws.Range('E48').Value = '  +12.31%  '
$ws.Range('D49').Value = '2.811.24'
$ws.Range('E49').Value = '  +0.82%  '
Set-TextValue 'D50' '103.96'
$ws.Range('E50').Value = '  +1.23%  '
$ws.Range('E51').Value = '  +3.82%  '
